$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - "Desenvolver o CRUD chamados": progress 0.75 -> 1
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = "READ,CREATE,UPDATE,DELETE"

# Row 10 - "Desenvolver o CRUD materiais": progress 0.75 -> 1
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = "READ,CREATE,UPDATE,DELETE"

# Row 11 - "Desenvolver o CRUD porta a porta": progress 0.75 -> 1
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = "READ,CREATE,UPDATE,DELETE"

# Row 12 - "Desenvolver o APP user comum": progress 0.5 -> 1
$ws.Range("N12").Value = 1

# Row 13 - "Desenvolver documentação": progress 0.5 -> 0.75, comment "EM ANDAMENTO" -> "FASE FINAL"
$ws.Range("N13").Value = 0.75
$ws.Range("O13").Value = "FASE FINAL"

# Row 8 - comment update (case/accent change)
$ws.Range("O8").Value = "TODOS TESTARÃO NO INSÔMINIA"

# Move the active selection cursor to N14 (cosmetic, matches source selection change)
[void]$ws.Range("N14").Select()
